$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update existing "Egypt industry-average" row, add two new growth columns ---
$ws.Cells.Item(2, 2).Value = "2"
$ws.Cells.Item(2, 4).Value = 0.2475
$ws.Cells.Item(2, 5).Value = 0.5945
$ws.Cells.Item(2, 7).Value = 0.2382513661202185
$ws.Cells.Item(2, 8).Value = 0.2382513661202185
$ws.Cells.Item(2, 9).Value = 0.2496174863387978
$ws.Cells.Item(2, 10).Value = 0.2416308869969101
$ws.Cells.Item(2, 11).Value = 21.99
$ws.Cells.Item(2, 12).Value = 0.240327868852459
$ws.Cells.Item(2, 13).Value = 5.84
$ws.Cells.Item(2, 14).Value = 0.04936601859678783
$ws.Cells.Item(2, 15).Value = 0.265575261482492
$ws.Cells.Item(2, 16).Value = 5.84
$ws.Cells.Item(2, 17).Value = 0.04936601859678783
$ws.Cells.Item(2, 18).Value = 0.265575261482492
$ws.Cells.Item(2, 21).Value = 13.28
$ws.Cells.Item(2, 22).Value = 0.1122569737954354
$ws.Cells.Item(2, 23).Value = 0.387280701754386
$ws.Cells.Item(2, 24).Value = 0.07572266396924258
$ws.Cells.Item(2, 25).Value = 0.3115580377851434
$ws.Cells.Item(2, 26).Value = 2.190567392865693
$ws.Cells.Item(2, 27).Value = 0.5870397756462727
$ws.Cells.Item(2, 28).Value = 0.07572266396924258
$ws.Cells.Item(2, 29).Value = 0.5113171116770302
$ws.Cells.Item(2, 33).Value = -13.28
$ws.Cells.Item(2, 36).Value = -0.1264521043610741
$ws.Cells.Item(2, 37).Value = -0.1955241460541814
$ws.Cells.Item(2, 38).Value = 0.078
$ws.Cells.Item(2, 39).Value = 0.078
$ws.Cells.Item(2, 41).Value = 292.8205128205128
$ws.Cells.Item(2, 42).Value = -0.5724137931034483
$ws.Cells.Item(2, 43).Value = 292.8205128205128

# --- Insert a new row at position 3; this shifts the old row 3 (Mohandes) down to row 4 ---
$ws.Rows.Item(3).Insert()

# --- Row 3 (new): Delta Insurance Company (CASE:DEIN) ---
$ws.Cells.Item(3, 1).Value = "Egypt"
$ws.Cells.Item(3, 2).Value = "Delta Insurance Company (CASE:DEIN)"
$ws.Cells.Item(3, 3).Value = "Insurance (General)"
$ws.Cells.Item(3, 4).Value = 0.323
$ws.Cells.Item(3, 5).Value = 0.919
$ws.Cells.Item(3, 7).Value = 0.2195121951219512
$ws.Cells.Item(3, 8).Value = 0.2195121951219512
$ws.Cells.Item(3, 9).Value = 0.2439024390243903
$ws.Cells.Item(3, 10).Value = 0.2381357534505578
$ws.Cells.Item(3, 11).Value = 12.6
$ws.Cells.Item(3, 12).Value = 0.2363977485928705
$ws.Cells.Item(3, 13).Value = -0
$ws.Cells.Item(3, 14).Value = -0
$ws.Cells.Item(3, 15).Value = -0
$ws.Cells.Item(3, 16).Value = -0
$ws.Cells.Item(3, 17).Value = -0
$ws.Cells.Item(3, 18).Value = -0
$ws.Cells.Item(3, 19).Value = 0
$ws.Cells.Item(3, 21).Value = 8.65
$ws.Cells.Item(3, 22).Value = 0.173
$ws.Cells.Item(3, 23).Value = 0.5
$ws.Cells.Item(3, 24).Value = 0.07572266396924258
$ws.Cells.Item(3, 25).Value = 0.4242773360307574
$ws.Cells.Item(3, 26).Value = 3.414477898782831
$ws.Cells.Item(3, 27).Value = 0.8131092670669269
$ws.Cells.Item(3, 28).Value = 0.07572266396924258
$ws.Cells.Item(3, 29).Value = 0.7373866030976843
$ws.Cells.Item(3, 30).Value = 0
$ws.Cells.Item(3, 31).Value = 0
$ws.Cells.Item(3, 32).Value = 0
$ws.Cells.Item(3, 33).Value = -8.65
$ws.Cells.Item(3, 34).Value = 0
$ws.Cells.Item(3, 35).Value = 0
$ws.Cells.Item(3, 36).Value = -0.2091898428053204
$ws.Cells.Item(3, 37).Value = -0.2562962962962963
$ws.Cells.Item(3, 38).Value = 0
$ws.Cells.Item(3, 39).Value = 0
$ws.Cells.Item(3, 40).Value = 0
$ws.Cells.Item(3, 42).Value = -0.6553030303030304

# --- Row 4: update shifted "Mohandes" row with refreshed figures ---
$ws.Cells.Item(4, 4).Value = 0.172
$ws.Cells.Item(4, 5).Value = 0.27
$ws.Cells.Item(4, 7).Value = 0.2643979057591623
$ws.Cells.Item(4, 8).Value = 0.2643979057591623
$ws.Cells.Item(4, 9).Value = 0.2575916230366492
$ws.Cells.Item(4, 10).Value = 0.247198498307387
$ws.Cells.Item(4, 11).Value = 9.390000000000001
$ws.Cells.Item(4, 12).Value = 0.2458115183246073
$ws.Cells.Item(4, 13).Value = 5.84
$ws.Cells.Item(4, 14).Value = 0.08550512445095168
$ws.Cells.Item(4, 15).Value = 0.6219382321618743
$ws.Cells.Item(4, 16).Value = 5.84
$ws.Cells.Item(4, 17).Value = 0.08550512445095168
$ws.Cells.Item(4, 18).Value = 0.6219382321618743
$ws.Cells.Item(4, 21).Value = 4.63
$ws.Cells.Item(4, 22).Value = 0.0677891654465593
$ws.Cells.Item(4, 23).Value = 0.274561403508772
$ws.Cells.Item(4, 24).Value = 0.07572266396924258
$ws.Cells.Item(4, 25).Value = 0.1988387395395294
$ws.Cells.Item(4, 26).Value = 1.460244648318043
$ws.Cells.Item(4, 27).Value = 0.3609702842256186
$ws.Cells.Item(4, 28).Value = 0.07572266396924258
$ws.Cells.Item(4, 29).Value = 0.285247620256376
$ws.Cells.Item(4, 33).Value = -4.63
$ws.Cells.Item(4, 36).Value = -0.07271870582692005
$ws.Cells.Item(4, 37).Value = -0.1354989757096869
$ws.Cells.Item(4, 38).Value = 0.078
$ws.Cells.Item(4, 39).Value = 0.078
$ws.Cells.Item(4, 41).Value = 126.1538461538461
$ws.Cells.Item(4, 42).Value = -0.463
$ws.Cells.Item(4, 43).Value = 126.1538461538461
